$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Mahle---Knecht / 02943N0 / MZYJ / 51 / 1510
$ws.Range("A2").Value = "Mahle---Knecht"
$ws.Range("B2").Value = "02943N0"
$ws.Range("C2").Value = "MZYJ"
$ws.Range("D2").Value = 51
$ws.Range("E2").Value = 1510

# Update row 3: Peugeot---Citroen / 82026 / JFWU / 4 / 309
$ws.Range("A3").Value = "Peugeot---Citroen"
$ws.Range("B3").Value = 82026
$ws.Range("C3").Value = "JFWU"
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 309

# Remove the now-unused rows 4-8
$ws.Range("A4:E8").Delete()
